# Update column F (dSF) values to match re-pulled / re-pushed data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    6  = -5
    7  = 1
    8  = -7
    9  = -1
    10 = -5
    11 = 5
    12 = -5
    14 = -8
    16 = -2
    17 = -3
    18 = -11
    19 = -6
    20 = 0
    21 = -11
    22 = -4
    24 = -3
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
